# Updates cryptos list prices / 1h volume percentages (and swaps the
# WrappedBTC / TRON rows 17 and 18) to match the latest scrape.
#
# Several "Price" column values look numeric (e.g. "1.00", "0.117") but
# must stay as plain text, matching the source workbook where every D/E
# cell is stored as a string. Setting NumberFormat to "@" (Text) right
# before assigning such values stops Excel from auto-converting them to
# numbers; resetting the Style back to "Normal" afterwards keeps the
# cell's style identical to the original (no explicit style index).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.642.27'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '3.519.30'
$ws.Range('E3').Value = '  -3.07%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.43%  '
$ws.Range('D7').Value = '3.518.32'
$ws.Range('E7').Value = '  -3.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.506'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.68'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('E11').Value = '  -5.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.404'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.25%  '
$ws.Range('D13').Value = '4.112.80'
$ws.Range('E13').Value = '  -3.12%  '
$ws.Range('E14').Value = '  -7.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.68'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.70%  '
$ws.Range('D16').Value = '3.531.16'
$ws.Range('E16').Value = '  -2.61%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '66.510.25'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.117'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -8.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '420.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.16%  '
$ws.Range('E23').Value = '  -5.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '76.86'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.54%  '
$ws.Range('D25').Value = '3.665.44'
$ws.Range('E25').Value = '  -2.88%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -8.42%  '
$ws.Range('E28').Value = '  -2.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.93'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').Value = '3.526.16'
$ws.Range('E32').Value = '  -2.82%  '
$ws.Range('E33').Value = '  -3.44%  '
$ws.Range('E36').Value = '  -10.35%  '
$ws.Range('E37').Value = '  -5.11%  '
$ws.Range('E38').Value = '  -5.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '172.97'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -9.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0803'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.18%  '
$ws.Range('E42').Value = '  -5.82%  '
$ws.Range('E43').Value = '  -5.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '45.58'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('E45').Value = '  -7.34%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.36'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.04'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('E49').Value = '  -5.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.85'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.96%  '
$ws.Range('E51').Value = '  -7.00%  '
